$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F header: "time_taken", styled like the other header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# Fill F2:F25 with the recorded time_taken timestamps (plain strings, not dates)
$ws.Range("F2").Value = "2021-10-05 10:52:31.057802"
$ws.Range("F3").Value = "2021-10-05 10:52:31.057813"
$ws.Range("F4").Value = "2021-10-05 10:52:31.057817"
$ws.Range("F5").Value = "2021-10-05 10:52:31.057819"
$ws.Range("F6").Value = "2021-10-05 10:52:31.057822"
$ws.Range("F7").Value = "2021-10-05 10:52:31.057825"
$ws.Range("F8").Value = "2021-10-05 10:52:31.057828"
$ws.Range("F9").Value = "2021-10-05 10:52:31.057831"
$ws.Range("F10").Value = "2021-10-05 10:52:31.057833"
$ws.Range("F11").Value = "2021-10-05 10:52:31.057836"
$ws.Range("F12").Value = "2021-10-05 10:52:31.057839"
$ws.Range("F13").Value = "2021-10-05 10:52:31.057841"
$ws.Range("F14").Value = "2021-10-05 10:52:31.057844"
$ws.Range("F15").Value = "2021-10-05 10:52:31.057846"
$ws.Range("F16").Value = "2021-10-05 10:52:31.057849"
$ws.Range("F17").Value = "2021-10-05 10:52:31.057852"
$ws.Range("F18").Value = "2021-10-05 10:52:31.057855"
$ws.Range("F19").Value = "2021-10-05 10:52:31.057857"
$ws.Range("F20").Value = "2021-10-05 10:52:31.057860"
$ws.Range("F21").Value = "2021-10-05 10:52:31.057862"
$ws.Range("F22").Value = "2021-10-05 10:52:31.057865"
$ws.Range("F23").Value = "2021-10-05 10:52:31.057867"
$ws.Range("F24").Value = "2021-10-05 10:52:31.057870"
$ws.Range("F25").Value = "2021-10-05 10:52:31.057873"

$excel.CutCopyMode = 0
